# Commit: "added fuel, adjusted import"
#
# Inserts a new "fuel" worksheet right after "asymmetric_sgen" (and right
# before "ext_grid"), listing the fuel type ("solar") for each of the 8
# sgen rows already present on the "sgen" sheet.

$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item("asymmetric_sgen")
$fuelSheet = $wb.Worksheets.Add($null, $afterSheet)
$fuelSheet.Name = "fuel"

# Fill the data rows first (column by column) so the shared-string table
# picks up "solar" / "sgen" before the header labels "fuel" / "index" /
# "gen_type" are written.
for ($i = 0; $i -le 7; $i++) {
    $row = $i + 2
    $fuelSheet.Cells.Item($row, 1).Value = $i
    $fuelSheet.Cells.Item($row, 4).Value = "solar"
}
for ($i = 0; $i -le 7; $i++) {
    $row = $i + 2
    $fuelSheet.Cells.Item($row, 2).Value = "sgen"
    $fuelSheet.Cells.Item($row, 3).Value = $i
}

# Header row, added last.
$fuelSheet.Range("D1").Value = "fuel"
$fuelSheet.Range("C1").Value = "index"
$fuelSheet.Range("B1").Value = "gen_type"
$fuelSheet.Range("B1:D1").Font.Bold = $true

# Match the saved selection/active cell on the new sheet.
[void]$fuelSheet.Range("K27").Select()

# The new sheet becomes the active tab of the workbook.
[void]$fuelSheet.Activate()
